$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append one new row (row 3) with the submitted trip record, mirroring the
# existing header/row-2 layout: Notes, Escort, Quantity, Camp, Trip type,
# Vehicle, Organization, Timestamp.
#
# A3 ("Notes") has no value for this submission, same as A2 above it - leave
# it blank, matching the other empty cell in that column.
$ws.Range("B3").Value = "أحمد شريم"

# The quantity column stores numbers-that-look-like-numbers as text (see the
# numberStoredAsText ignoredError covering the whole sheet range), so force
# a text number format before writing it - otherwise "500" would be
# reinterpreted as a numeric value instead of staying text like "2" in C2.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "500"

$ws.Range("D3").Value = "الصمود"
$ws.Range("E3").Value = "الرحلة 3"
$ws.Range("F3").Value = "C5"
$ws.Range("G3").Value = "UNICEF"
$ws.Range("H3").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٣٨:٥٥ م"
